$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match header formatting used by existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for columns I and J, rows 2-16
$data = @{
    2  = @(1, 6)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 4)
    6  = @(1, 4)
    7  = @(3, 4)
    8  = @(3, 6)
    9  = @(1, 3)
    10 = @(1, 3)
    11 = @(1, 3)
    12 = @(7, 8)
    13 = @(2, 6)
    14 = @(1, 4)
    15 = @(1, 3)
    16 = @(8, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
